$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated TPM-derived values (ligand avg/total expression, receptor avg/total
# expression, and the dependent specificity / edge-weight columns) for rows 2-7.

$ws.Range("G2").Value = 0.159713
$ws.Range("H2").Value = 0.319426
$ws.Range("M2").Value = 4.296436999999999
$ws.Range("N2").Value = 8.592873999999998
$ws.Range("O2").Value = 0.08737129157293111
$ws.Range("P2").Value = 0.06876644796033347
$ws.Range("Q2").Value = 0.6861968425809999
$ws.Range("R2").Value = 2.744787370324
$ws.Range("S2").Value = 0.08737129157293111
$ws.Range("T2").Value = 0.06876644796033347

$ws.Range("G3").Value = 0.159713
$ws.Range("H3").Value = 0.319426
$ws.Range("O3").Value = 0.2709785829485105
$ws.Range("P3").Value = 0.3199146015909443
$ws.Range("Q3").Value = 2.128212192801667
$ws.Range("R3").Value = 12.76927315681
$ws.Range("S3").Value = 0.2709785829485105
$ws.Range("T3").Value = 0.3199146015909443

$ws.Range("G4").Value = 0.159713
$ws.Range("H4").Value = 0.319426
$ws.Range("M4").Value = 6.89049
$ws.Range("N4").Value = 20.67147
$ws.Range("O4").Value = 0.140123318663899
$ws.Range("P4").Value = 0.1654281868928364
$ws.Range("Q4").Value = 1.10050082937
$ws.Range("R4").Value = 6.603004976219999
$ws.Range("S4").Value = 0.140123318663899
$ws.Range("T4").Value = 0.1654281868928364

$ws.Range("G5").Value = 0.159713
$ws.Range("H5").Value = 0.319426
$ws.Range("M5").Value = 18.2696115
$ws.Range("N5").Value = 36.539223
$ws.Range("O5").Value = 0.3715263492262718
$ws.Range("P5").Value = 0.292413525083752
$ws.Range("Q5").Value = 2.9178944614995
$ws.Range("R5").Value = 11.671577845998
$ws.Range("S5").Value = 0.3715263492262718
$ws.Range("T5").Value = 0.292413525083752

$ws.Range("G6").Value = 0.159713
$ws.Range("H6").Value = 0.319426
$ws.Range("M6").Value = 1.355562
$ws.Range("N6").Value = 4.066686000000001
$ws.Range("O6").Value = 0.02756637715092428
$ws.Range("P6").Value = 0.03254458882907125
$ws.Range("Q6").Value = 0.216500873706
$ws.Range("R6").Value = 1.299005242236
$ws.Range("S6").Value = 0.02756637715092428
$ws.Range("T6").Value = 0.03254458882907125

$ws.Range("G7").Value = 0.159713
$ws.Range("H7").Value = 0.319426
$ws.Range("M7").Value = 5.037141666666667
$ws.Range("N7").Value = 15.111425
$ws.Range("O7").Value = 0.1024340804374633
$ws.Range("P7").Value = 0.1209326496430627
$ws.Range("Q7").Value = 0.8044970070083334
$ws.Range("R7").Value = 4.82698204205
$ws.Range("S7").Value = 0.1024340804374633
$ws.Range("T7").Value = 0.1209326496430627
